# Trade #62 closed at 2026-02-17 12:53:11 - unknown UNKNOWN +0.000%
#
# 1) Summary sheet: bump Total Trades (B6) and recompute Win Rate % (B9)
# 2) Strategy Status sheet: same bump for the MarketMaking row (D4 / G4)
# 3) All Trades + MarketMaking sheets: append the new trade as row 63
#
# Note on the new row's Date/Time cells (B63/C63): assigning a literal
# string like "2026-02-17" straight to .Value makes Excel's type-inference
# treat it as a real date/time serial (and stamps a non-default
# NumberFormat style on the cell) - the existing rows store these as plain
# text instead. To avoid that:
#  - B63 (Date) is identical to every other row's date ("2026-02-17"), so
#    we just .Copy the cell above it straight down - value and (default)
#    format, no re-parsing involved.
#  - C63 (Time) is a new value, so we first .Copy the cell above it (which
#    seeds C63 as plain text with the default style) and only then
#    overwrite .Value - Excel keeps treating an already-text cell's new
#    value as text instead of re-parsing a HH:MM:SS-shaped string as a
#    time, the way it would for a still-empty/default cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 62
$summary.Range("B9").Value = 43.55

# ---------------------------------------------------------------------
# Strategy Status (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 62
$status.Range("G4").Value = 43.55

# ---------------------------------------------------------------------
# Helper: append trade #62 as row 63 on a trades-log sheet
# ---------------------------------------------------------------------
function Add-Trade62Row($ws) {
    # Date is unchanged from the row above - plain copy keeps it text.
    $ws.Range("B62").Copy($ws.Range("B63"))

    # Time is new - seed as text via Copy, then overwrite the value so it
    # does not get re-parsed into a time serial.
    $ws.Range("C62").Copy($ws.Range("C63"))
    $ws.Range("C63").Value = "12:53:05"

    $ws.Range("A63").Value = 62
    $ws.Range("D63").Value = "MarketMaking"
    $ws.Range("E63").Value = "DOWN"
    $ws.Range("F63").Value = 0.86
    $ws.Range("G63").Value = 0.86
    $ws.Range("H63").Value = "CLOSED"
    $ws.Range("I63").Value = 0
    $ws.Range("J63").Value = 0
    $ws.Range("K63").Value = 100.11
    $ws.Range("L63").Value = 0
    $ws.Range("M63").Value = 0
    $ws.Range("N63").Value = 0.6
    $ws.Range("O63").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P63").Value = "early_exit"
    $ws.Range("Q63").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade62Row $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade62Row $marketMaking
